# Update Breakdowns_List worksheet with new breakdown timestamps (including
# full days, not just intra-day corrections) and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 43795.447222222225
$ws.Range("B2").Value = 43795.602777777778
$ws.Range("C2").Value = 43807.884722222225

# Row 3
$ws.Range("A3").Value = 43796.511805555558
$ws.Range("B3").Value = 43796.568749999999
$ws.Range("C3").Value = 43799.341666666667

# Row 4
$ws.Range("A4").Value = 43795.583333333336
$ws.Range("B4").Value = 43796.333333333336
$ws.Range("C4").Value = 43827.916666666664

# Update the active selection to D4 (was E5)
$ws.Range("D4").Select()
